$wb = $excel.ActiveWorkbook

$wsDetectors = $wb.Worksheets.Item("Detectors")
$wsTriggers  = $wb.Worksheets.Item("Triggers")

# --- Fix the "Op Effiency (J/op)" typo on the Detectors sheet header.
# Overwriting the cell drops the last reference to the mis-spelled shared
# string so it gets pruned from the shared-string table on save, and every
# other cell that pointed at the correctly spelled string is renumbered
# automatically (Triggers!H1, Global!A1).
$wsDetectors.Range("G1").Value = "Op Efficiency (J/op)"

# --- Triggers sheet: swap the "Output" / "Name" columns (A <-> B) ---------
$wsTriggers.Range("A1").Value = "Name"
$wsTriggers.Range("B1").Value = "Output"

$wsTriggers.Range("A2").Value = "Tracking"
$wsTriggers.Range("B2").Value = "Intermediate"

$wsTriggers.Range("A3").Value = "Timing"
$wsTriggers.Range("B3").Value = "Intermediate"

$wsTriggers.Range("A4").Value = "Calorimetry"
$wsTriggers.Range("B4").Value = "Intermediate"

$wsTriggers.Range("A5").Value = "Muon"
$wsTriggers.Range("B5").Value = "Intermediate"

$wsTriggers.Range("A6").Value = "Intermediate"
$wsTriggers.Range("B6").Value = "Global"

$wsTriggers.Range("A7").Value = "Global"
$wsTriggers.Range("B7").Value = "Disk"

$wsTriggers.Range("A8").Value = "Disk"
$wsTriggers.Range("B8").Value = "None"

# --- Triggers sheet: power-calculation refactor values ---------------------
$wsTriggers.Range("H6").Value = 0.003
$wsTriggers.Range("H7").Value = 16

# --- Triggers sheet: formatting -------------------------------------------
# Whole table gets the explicit-black-font style (fontId 1), and column G
# keeps/gets the scientific number format on top of that font.
$wsTriggers.Range("A1:I8").Font.Color = 0
$wsTriggers.Range("G2:G8").NumberFormat = "0.00E+00"

# --- Sheet view / selection bookkeeping ------------------------------------
$wsTriggers.Activate()
$wsTriggers.Range("E20").Select()

$wsDetectors.Activate()
$wsDetectors.Range("J17").Select()
